# Fix Training Data Issue (#48)
# The BF column ("Date") holds the game date as free text, e.g. "6-28-2007-08".
# Because of how the NBA stats site showed the season label, the date was
# off by one day / formatted oddly. Correct it to ISO-ish "2008-06-28" text
# for every data row (rows 2-31; row 1 is the "Date" header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 31
$col      = "BF"

# Pre-format the target cells as Text so Excel's auto-detection doesn't
# reinterpret the new value ("2008-06-28") as a date serial number.
$rng = $ws.Range("$col$firstRow`:$col$lastRow")
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("$col$r").Value = "2008-06-28"
}

# Restore the default (unformatted) style so the cells keep behaving like
# the rest of the plain, unstyled data cells.
$rng.Style = "Normal"
